# Rename the "_old"/"_new" header-suffix columns to "_FV2404"/"_FV2410"
# (format-version-specific) suffixes, then wrap the data range in a table
# and freeze the header row, matching the target workbook layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10) carry the "old" (FV2404) values, column K (11) is the
# "diff" column, and columns L-U (12-21) carry the "new" (FV2410) values.
for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $baseHeaders[$i] + "_FV2404"
}
for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, 11 + $i + 1).Value2 = $baseHeaders[$i] + "_FV2410"
}

# Turn the whole sheet's data range into an actual Excel table.
$dataRange = $ws.Range("A1:U69")
$listObject = $ws.ListObjects.Add(1, $dataRange, [Type]::Missing, 1, [Type]::Missing)
$listObject.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
